$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 28
$ws.Range("H28").Value = 6152.905
$ws.Range("I28").Value = 400.3846
$ws.Range("K28").Value = 400.3846
$ws.Range("M28").Value = 84.61540000000002

# row 32
$ws.Range("H32").Value = 1186.1765
$ws.Range("J32").Value = 1036.6923
$ws.Range("L32").Value = 1036.6923
$ws.Range("N32").Value = -1688.6923

# row 62
$ws.Range("H62").Value = 3250
$ws.Range("I62").Value = 1675
$ws.Range("J62").Value = 3880
$ws.Range("K62").Value = 1675
$ws.Range("L62").Value = 3880
$ws.Range("M62").Value = -1051
$ws.Range("N62").Value = -5128

# row 65
$ws.Range("H65").Value = 3250
$ws.Range("I65").Value = 1675
$ws.Range("J65").Value = 3880
$ws.Range("K65").Value = 8375
$ws.Range("L65").Value = 19400
$ws.Range("M65").Value = -5255
$ws.Range("N65").Value = -25640

# row 113
$ws.Range("H113").Value = 2573.182
$ws.Range("I113").Value = 2334.1667
$ws.Range("K113").Value = 2334.1667
$ws.Range("M113").Value = 919.8332999999998

# row 116
$ws.Range("H116").Value = 2630.7778
$ws.Range("I116").Value = 2507.5
$ws.Range("J116").Value = 2877.3333
$ws.Range("K116").Value = 2507.5
$ws.Range("L116").Value = 2877.3333
$ws.Range("M116").Value = 934.5
$ws.Range("N116").Value = -9761.3333

# row 132
$ws.Range("H132").Value = 4308.643
$ws.Range("I132").Value = 4315.1055
$ws.Range("J132").Value = 4295
$ws.Range("K132").Value = 12945.3165
$ws.Range("L132").Value = 12885
$ws.Range("M132").Value = -10415.3165
$ws.Range("N132").Value = -17945

$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 1164.0769
$ws.Range("I2").Value = 892.55554
$ws.Range("J2").Value = 1775
$ws.Range("K2").Value = 892.55554
$ws.Range("L2").Value = 1775
$ws.Range("M2").Value = -779.55554
$ws.Range("N2").Value = -2001

# row 45
$ws.Range("H45").Value = 3067.1667
$ws.Range("I45").Value = 2308
$ws.Range("J45").Value = 4130
$ws.Range("K45").Value = 2308
$ws.Range("L45").Value = 4130
$ws.Range("M45").Value = -1931
$ws.Range("N45").Value = -4884

# row 102
$ws.Range("H102").Value = 1750
$ws.Range("I102").Value = 1750
$ws.Range("K102").Value = 1750
$ws.Range("M102").Value = -128

# row 116
$ws.Range("H116").Value = 1164.0769
$ws.Range("I116").Value = 892.55554
$ws.Range("J116").Value = 1775
$ws.Range("K116").Value = 892.55554
$ws.Range("L116").Value = 1775
$ws.Range("M116").Value = 1401.44446
$ws.Range("N116").Value = -6363

$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 1164.0769
$ws.Range("I3").Value = 892.55554
$ws.Range("J3").Value = 1775
$ws.Range("K3").Value = 892.55554
$ws.Range("L3").Value = 1775
$ws.Range("M3").Value = -778.55554
$ws.Range("N3").Value = -2003

# row 80
$ws.Range("H80").Value = 1609.8462
$ws.Range("I80").Value = 4569
$ws.Range("J80").Value = 294.66666
$ws.Range("K80").Value = 4569
$ws.Range("L80").Value = 294.66666
$ws.Range("M80").Value = -3571
$ws.Range("N80").Value = -2290.66666

# row 83
$ws.Range("H83").Value = 1609.8462
$ws.Range("I83").Value = 4569
$ws.Range("J83").Value = 294.66666
$ws.Range("K83").Value = 22845
$ws.Range("L83").Value = 1473.3333
$ws.Range("M83").Value = -17853
$ws.Range("N83").Value = -11457.3333

# row 94
$ws.Range("H94").Value = 785.125
$ws.Range("I94").Value = 689.36365
$ws.Range("K94").Value = 689.36365
$ws.Range("M94").Value = -238.36365

# row 105
$ws.Range("H105").Value = 7814825.5
$ws.Range("J105").Value = 3075
$ws.Range("L105").Value = 3075
$ws.Range("N105").Value = -6569

# row 134
$ws.Range("H134").Value = 3834.7334
$ws.Range("I134").Value = 4253.5
$ws.Range("J134").Value = 3555.5557
$ws.Range("K134").Value = 12760.5
$ws.Range("L134").Value = 10666.6671
$ws.Range("M134").Value = -10225.5
$ws.Range("N134").Value = -15736.6671

$ws = $wb.Worksheets.Item("CRP")
# row 16
$ws.Range("H16").Value = 910.1429000000001
$ws.Range("I16").Value = 842.75
$ws.Range("K16").Value = 842.75
$ws.Range("M16").Value = -555.75

# row 113
$ws.Range("H113").Value = 910.1429000000001
$ws.Range("I113").Value = 842.75
$ws.Range("K113").Value = 842.75
$ws.Range("M113").Value = 1327.25

$ws = $wb.Worksheets.Item("CUL")
# row 35
$ws.Range("H35").Value = 4481.2
$ws.Range("I35").Value = 300
$ws.Range("J35").Value = 5526.5
$ws.Range("K35").Value = 900
$ws.Range("L35").Value = 16579.5
$ws.Range("M35").Value = -612
$ws.Range("N35").Value = -17155.5

# row 136
$ws.Range("H136").Value = 1887.7778
$ws.Range("I136").Value = 1455.7142
$ws.Range("J136").Value = 3400
$ws.Range("K136").Value = 4367.142599999999
$ws.Range("L136").Value = 10200
$ws.Range("M136").Value = 732.8574000000008
$ws.Range("N136").Value = -20400

$ws = $wb.Worksheets.Item("GSM")
# row 97
$ws.Range("H97").Value = 3017.2856
$ws.Range("I97").Value = 2977.5
$ws.Range("K97").Value = 2977.5
$ws.Range("M97").Value = -2481.5

# row 113
$ws.Range("H113").Value = 1335.5714
$ws.Range("I113").Value = 869.8
$ws.Range("K113").Value = 869.8
$ws.Range("M113").Value = 1300.2

# row 122
$ws.Range("H122").Value = 4901.143
$ws.Range("I122").Value = 3650.2
$ws.Range("J122").Value = 5596.1113
$ws.Range("K122").Value = 10950.6
$ws.Range("L122").Value = 16788.3339
$ws.Range("M122").Value = -8500.599999999999
$ws.Range("N122").Value = -21688.3339

$ws = $wb.Worksheets.Item("LTW")
# row 61
$ws.Range("H61").Value = 3470.7856
$ws.Range("I61").Value = 3511.4211
$ws.Range("J61").Value = 3385
$ws.Range("K61").Value = 3511.4211
$ws.Range("L61").Value = 3385
$ws.Range("M61").Value = -3309.4211
$ws.Range("N61").Value = -3789

# row 113
$ws.Range("H113").Value = 3470.7856
$ws.Range("I113").Value = 3511.4211
$ws.Range("J113").Value = 3385
$ws.Range("K113").Value = 3511.4211
$ws.Range("L113").Value = 3385
$ws.Range("M113").Value = -1341.4211
$ws.Range("N113").Value = -7725

# row 136
$ws.Range("H136").Value = 33336334
$ws.Range("I136").Value = 3750
$ws.Range("K136").Value = 11250
$ws.Range("M136").Value = -8700

$ws = $wb.Worksheets.Item("WVR")
# row 107
$ws.Range("H107").Value = 570.1111
$ws.Range("I107").Value = 648.7143
$ws.Range("J107").Value = 295
$ws.Range("K107").Value = 1946.1429
$ws.Range("L107").Value = 885
$ws.Range("M107").Value = -26.14289999999983
$ws.Range("N107").Value = -4725
